$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 121, shifting existing rows 121-135 down to 124-138.
$ws.Rows.Item(121).Resize(3).Insert()

# --- Row 121: Primera, Provincia de Melipilla, $/bandeja 7 kilos ---
$ws.Range("A121").Value = 10
$ws.Range("B121").Value = "Vega Modelo de Temuco"
$ws.Range("C121").Value = "La Araucanía"
$ws.Range("D121").Value = 44476
$ws.Range("E121").Value = 9
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100101
$ws.Range("H121").Value = "Berries"
$ws.Range("I121").Value = 100112025
$ws.Range("J121").Value = "Frutilla"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 4000
$ws.Range("N121").Value = 15000
$ws.Range("O121").Value = 15000
$ws.Range("P121").Value = 15000
$ws.Range("Q121").Value = "$/bandeja 7 kilos"
$ws.Range("R121").Value = "Provincia de Melipilla"
$ws.Range("S121").Value = 2143
$ws.Range("T121").Value = 7

# --- Row 122: Segunda, Provincia de Melipilla, $/bandeja 7 kilos ---
$ws.Range("A122").Value = 10
$ws.Range("B122").Value = "Vega Modelo de Temuco"
$ws.Range("C122").Value = "La Araucanía"
$ws.Range("D122").Value = 44476
$ws.Range("E122").Value = 9
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100101
$ws.Range("H122").Value = "Berries"
$ws.Range("I122").Value = 100112025
$ws.Range("J122").Value = "Frutilla"
$ws.Range("K122").Value = "Sin especificar"
$ws.Range("L122").Value = "Segunda"
$ws.Range("M122").Value = 400
$ws.Range("N122").Value = 13000
$ws.Range("O122").Value = 13000
$ws.Range("P122").Value = 13000
$ws.Range("Q122").Value = "$/bandeja 7 kilos"
$ws.Range("R122").Value = "Provincia de Melipilla"
$ws.Range("S122").Value = 1857
$ws.Range("T122").Value = 7

# --- Row 123: Tercera, Provincia de Melipilla, $/bandeja 7 kilos ---
$ws.Range("A123").Value = 10
$ws.Range("B123").Value = "Vega Modelo de Temuco"
$ws.Range("C123").Value = "La Araucanía"
$ws.Range("D123").Value = 44476
$ws.Range("E123").Value = 9
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100101
$ws.Range("H123").Value = "Berries"
$ws.Range("I123").Value = 100112025
$ws.Range("J123").Value = "Frutilla"
$ws.Range("K123").Value = "Sin especificar"
$ws.Range("L123").Value = "Tercera"
$ws.Range("M123").Value = 100
$ws.Range("N123").Value = 8000
$ws.Range("O123").Value = 8000
$ws.Range("P123").Value = 8000
$ws.Range("Q123").Value = "$/bandeja 7 kilos"
$ws.Range("R123").Value = "Provincia de Melipilla"
$ws.Range("S123").Value = 1143
$ws.Range("T123").Value = 7
